$d = $word.ActiveDocument

# Locate the paragraph that holds the "${invoices}" merge-field placeholder.
$targetIndex = $null
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like '*${invoices}*') {
        $targetIndex = $i
    }
}

if ($targetIndex -eq $null) {
    throw 'Could not find the paragraph containing the ${invoices} placeholder'
}

# The placeholder paragraph is followed by two empty "Title" styled
# paragraphs that together form the little block that will be expanded.
$firstPara = $d.Paragraphs.Item($targetIndex)
$lastPara  = $d.Paragraphs.Item($targetIndex + 2)

$blockRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$emptyPara = '<w:p ' + $wNs + '><w:pPr><w:jc w:val="center"/></w:pPr></w:p>'

$finalPara = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="25"/></w:rPr></w:pPr>' + `
    '<w:r><w:t>${invoice</w:t></w:r>' + `
    '<w:r><w:t>s</w:t></w:r>' + `
    '<w:r><w:t>}</w:t></w:r>' + `
    '</w:p>'

$newXml = ''
for ($n = 0; $n -lt 12; $n++) {
    $newXml = $newXml + $emptyPara
}
$newXml = $newXml + $finalPara

$null = $blockRange.InsertXML($newXml)
